$p = $ppt.ActivePresentation

$newText = "US 2: As an experienced player I want to have a greater diversity of towers that can help me automate the way I defend orders."

# The "US 2" user-story card is shape Id=5 on the Sprint board slides
# (slides 3 through 7 -- sldId 271, 278, 279, 280, 281). Locate it by its
# stable shape Id (not positional index, since the shapes were reordered)
# and correct the wording of the story.
for ($slideIndex = 3; $slideIndex -le 7; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.Id -eq 5) {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}
